# Add a new "M_HRD" column (C) with its header and data, and correct the
# Cen X-3 M_exp value which was fixed at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing value for Cen X-3 (row 7) in column B.
$ws.Range("B7").Value = 39

# New column header.
$ws.Range("C1").Value = "M_HRD"

# New column C values (M_HRD), rows 2-9.
$ws.Range("C2").Value = 33
$ws.Range("C3").Value = 36
$ws.Range("C4").Value = 21
$ws.Range("C5").Value = 45
$ws.Range("C6").Value = 19
$ws.Range("C7").Value = 48
$ws.Range("C8").Value = 23
$ws.Range("C9").Value = 12

# Update the selection to match where the author ended up (C9).
$ws.Range("C9").Select()

# Configure the page setup used for printing, as recorded in the saved file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
